$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = "46.313.87"
$ws.Cells.Item(2, 5).Value = "  +1.89%  "

$ws.Cells.Item(3, 4).Value = "2.614.71"
$ws.Cells.Item(3, 5).Value = "  +10.03%  "

$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  +0.05%  "

$ws.Cells.Item(5, 4).Value = "'313.25"
$ws.Cells.Item(5, 5).Value = "  +4.82%  "

$ws.Cells.Item(6, 4).Value = "'100.64"
$ws.Cells.Item(6, 5).Value = "  +4.22%  "

$ws.Cells.Item(7, 4).Value = "'0.598"
$ws.Cells.Item(7, 5).Value = "  +6.97%  "

$ws.Cells.Item(8, 5).Value = "  +0.12%  "

$ws.Cells.Item(9, 4).Value = "'0.584"
$ws.Cells.Item(9, 5).Value = "  +16.66%  "

$ws.Cells.Item(10, 4).Value = "'38.76"
$ws.Cells.Item(10, 5).Value = "  +14.06%  "

$ws.Cells.Item(11, 4).Value = "'0.0842"
$ws.Cells.Item(11, 5).Value = "  +7.48%  "

$ws.Cells.Item(12, 4).Value = "'8.36"
$ws.Cells.Item(12, 5).Value = "  +18.93%  "

$ws.Cells.Item(13, 4).Value = "3.013.90"
$ws.Cells.Item(13, 5).Value = "  +10.01%  "

$ws.Cells.Item(14, 5).Value = "  +1.92%  "

$ws.Cells.Item(15, 4).Value = "2.619.29"
$ws.Cells.Item(15, 5).Value = "  +10.06%  "

$ws.Cells.Item(16, 5).Value = "  +11.56%  "

$ws.Cells.Item(17, 4).Value = "'15.01"
$ws.Cells.Item(17, 5).Value = "  +8.95%  "

$ws.Cells.Item(18, 4).Value = "46.519.42"
$ws.Cells.Item(18, 5).Value = "  +2.37%  "

$ws.Cells.Item(19, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(19, 4).Value = "'13.24"
$ws.Cells.Item(19, 5).Value = "  +4.23%  "

$ws.Cells.Item(20, 2).Value = "ShibaInu"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(20, 4).Value = "'0.0000102"
$ws.Cells.Item(20, 5).Value = "  +8.47%  "

$ws.Cells.Item(21, 4).Value = "'6.75"
$ws.Cells.Item(21, 5).Value = "  +11.88%  "

$ws.Cells.Item(22, 4).Value = "'71.30"
$ws.Cells.Item(22, 5).Value = "  +6.79%  "

$ws.Cells.Item(23, 4).Value = "'255.82"
$ws.Cells.Item(23, 5).Value = "  +6.21%  "

$ws.Cells.Item(24, 4).Value = "'3.09"
$ws.Cells.Item(24, 5).Value = "  +12.39%  "

$ws.Cells.Item(25, 5).Value = "  +17.55%  "

$ws.Cells.Item(26, 4).Value = "'28.27"
$ws.Cells.Item(26, 5).Value = "  +36.24%  "

$ws.Cells.Item(27, 5).Value = "  -0.03%  "

$ws.Cells.Item(28, 4).Value = "'10.62"
$ws.Cells.Item(28, 5).Value = "  +10.32%  "

$ws.Cells.Item(29, 4).Value = "'40.16"
$ws.Cells.Item(29, 5).Value = "  +3.85%  "

$ws.Cells.Item(30, 4).Value = "'2.26"
$ws.Cells.Item(30, 5).Value = "  +2.68%  "

$ws.Cells.Item(31, 2).Value = "Filecoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(31, 4).Value = "'6.18"
$ws.Cells.Item(31, 5).Value = "  +13.29%  "

$ws.Cells.Item(32, 2).Value = "LidoDAOToken"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(32, 4).Value = "'3.73"
$ws.Cells.Item(32, 5).Value = "  -1.49%  "

$ws.Cells.Item(33, 5).Value = "  +20.69%  "

$ws.Cells.Item(34, 5).Value = "  +6.52%  "

$ws.Cells.Item(35, 4).Value = "'153.17"
$ws.Cells.Item(35, 5).Value = "  +4.17%  "

$ws.Cells.Item(36, 4).Value = "'0.0837"
$ws.Cells.Item(36, 5).Value = "  +9.46%  "

$ws.Cells.Item(37, 5).Value = "  +6.01%  "

$ws.Cells.Item(38, 5).Value = "  +5.94%  "

$ws.Cells.Item(39, 4).Value = "'17.18"
$ws.Cells.Item(39, 5).Value = "  +13.79%  "

$ws.Cells.Item(40, 4).Value = "'4.22"
$ws.Cells.Item(40, 5).Value = "  +10.55%  "

$ws.Cells.Item(41, 5).Value = "  +13.13%  "

$ws.Cells.Item(42, 4).Value = "'0.0326"
$ws.Cells.Item(42, 5).Value = "  +10.49%  "

$ws.Cells.Item(43, 4).Value = "'21.10"
$ws.Cells.Item(43, 5).Value = "  +48.84%  "

$ws.Cells.Item(44, 4).Value = "2.046.65"
$ws.Cells.Item(44, 5).Value = "  +5.14%  "

$ws.Cells.Item(45, 5).Value = "  +0.09%  "

$ws.Cells.Item(46, 4).Value = "'91.46"
$ws.Cells.Item(46, 5).Value = "  -1.07%  "

$ws.Cells.Item(47, 4).Value = "'9.30"
$ws.Cells.Item(47, 5).Value = "  +8.96%  "

$ws.Cells.Item(48, 2).Value = "Stacks"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(48, 4).Value = "'1.80"
$ws.Cells.Item(48, 5).Value = "  +2.50%  "

$ws.Cells.Item(49, 2).Value = "Aave"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(49, 4).Value = "'110.43"
$ws.Cells.Item(49, 5).Value = "  +12.60%  "

$ws.Cells.Item(50, 2).Value = "Algorand"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(50, 4).Value = "'0.202"
$ws.Cells.Item(50, 5).Value = "  +10.87%  "

$ws.Cells.Item(51, 2).Value = "RocketPoolETH"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(51, 4).Value = "2.873.79"
$ws.Cells.Item(51, 5).Value = "  +10.14%  "
